$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the last existing row (40) into the new row (41), which
# brings along matching cell styles (date/time formats, alignment, etc.)
# and shifts nothing else down since row 41 was previously blank.
$ws.Rows.Item(40).Copy()
$ws.Rows.Item(41).Insert(-4121)   # xlShiftDown

# Fill in the new entry: VGA Top / Arch, 11:00 - 11:30, Add ROM 1
$ws.Range("B41").Value2 = 0.45833333333333331   # 11:00
$ws.Range("C41").Value2 = 0.47916666666666669   # 11:30
$ws.Range("G41").Value2 = "Add ROM 1"

# Move the active selection to the newly added note cell
[void]$ws.Range("G41").Select()
